$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet 1: 展览  (numeric "想去人数" / interest-count updates only)
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("F2").Value  = 2865
$ws1.Range("F3").Value  = 1162
$ws1.Range("F4").Value  = 20973
$ws1.Range("F6").Value  = 2776
$ws1.Range("F7").Value  = 799
$ws1.Range("F9").Value  = 508
$ws1.Range("F10").Value = 760
$ws1.Range("F15").Value = 510
$ws1.Range("F18").Value = 14
$ws1.Range("F19").Value = 418
$ws1.Range("F20").Value = 49
$ws1.Range("F23").Value = 25

# ---------------------------------------------------------------------------
# Sheet 2: 演出  (numeric "想去人数" updates only)
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("F5").Value  = 333
$ws2.Range("F10").Value = 15
$ws2.Range("F14").Value = 146
$ws2.Range("F22").Value = 39

# ---------------------------------------------------------------------------
# Sheet 3: 本地生活  (numeric "想去人数" updates only)
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item(3)
$ws3.Range("F2").Value = 6119
$ws3.Range("F3").Value = 695
$ws3.Range("F4").Value = 679
$ws3.Range("F5").Value = 1557

# ---------------------------------------------------------------------------
# Sheet 4: 全部类型  (numeric "想去人数" updates ...)
# ---------------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("F2").Value  = 6119
$ws4.Range("F3").Value  = 695
$ws4.Range("F4").Value  = 679
$ws4.Range("F5").Value  = 1557
$ws4.Range("F6").Value  = 2865
$ws4.Range("F7").Value  = 1162
$ws4.Range("F8").Value  = 20973
$ws4.Range("F13").Value = 333
$ws4.Range("F14").Value = 2776
$ws4.Range("F15").Value = 799
$ws4.Range("F19").Value = 508
$ws4.Range("F20").Value = 760
$ws4.Range("F28").Value = 15
$ws4.Range("F30").Value = 510
$ws4.Range("F44").Value = 25
$ws4.Range("F49").Value = 39

# ... plus a block of rows 33-39 whose event listings were refreshed: the
# expired "KOKIA" event (row 33) dropped off, every following row shifted
# up by one, and a brand new event ("环形宇宙动漫游戏嘉年华") was appended
# as the new row 39.

# Row 33 <- old row 34 content (原神x星穹铁道x绝区零同人ONLY)
$ws4.Range("B33").Value = "'2024-11-09"
$ws4.Range("C33").Value = "广州·原神×星穹铁道×绝区零同人ONLY"
$ws4.Range("D33").Value = "西环路1号 广州岭南会展中心"
$ws4.Range("E33").Value = "2024.11.09 10:00-11.09 17:00"
$ws4.Range("F33").Value = 257
$ws4.Range("G33").Value = 60
$ws4.Range("H33").Value = "https://show.bilibili.com/platform/detail.html?id=92397"
$ws4.Range("I33").Value = "//i1.hdslb.com/bfs/openplatform/202409/t5ilbPxZ1726023971484.jpeg"

# Row 34 <- old row 35 content (平田雄也&小池亮介2024粉丝见面会), F bumped to 146
$ws4.Range("B34").Value = "'2024-11-10"
$ws4.Range("C34").Value = "广州·平田雄也&小池亮介2024粉丝见面会"
$ws4.Range("D34").Value = "金花街道中山七路333号1906科技圆区3号楼109-1铺、110-1铺、111-1铺 音乐唐人馆"
$ws4.Range("E34").Value = "2024.11.10 13:00-11.10 18:00"
$ws4.Range("F34").Value = 146
$ws4.Range("G34").Value = 480
$ws4.Range("H34").Value = "https://show.bilibili.com/platform/detail.html?id=92655"
$ws4.Range("I34").Value = "//i2.hdslb.com/bfs/openplatform/202409/UkhOeOwe1726658317935.jpeg"

# Row 35: same event as row 34 (duplicate listing kept as-is), only F changes
$ws4.Range("F35").Value = 146

# Row 36 <- old row 37 content (BanG Dream!only2·浮想宣告)
$ws4.Range("B36").Value = "'2024-11-16"
$ws4.Range("C36").Value = "广州·BanG Dream!only2·浮想宣告"
$ws4.Range("D36").Value = "同泰路颐和山庄内 颐和山庄国际会议厅"
$ws4.Range("E36").Value = "2024.11.16 10:00-11.16 17:00"
$ws4.Range("F36").Value = 14
$ws4.Range("G36").Value = 89
$ws4.Range("H36").Value = "https://show.bilibili.com/platform/detail.html?id=93056"
$ws4.Range("I36").Value = "//i1.hdslb.com/bfs/openplatform/202409/PaoNiZxp1727595871985.jpeg"

# Row 37 <- old row 38 content (wio jumponly4.0万圣狂欢节); B37 date unchanged
$ws4.Range("C37").Value = "广州·wio jumponly4.0万圣狂欢节"
$ws4.Range("D37").Value = "逸景路462号珠江国际纺织城d区6层 珠江时尚馆"
$ws4.Range("E37").Value = "2024.11.16 10:00-11.17 17:00"
$ws4.Range("F37").Value = 418
$ws4.Range("G37").Value = 23.3
$ws4.Range("H37").Value = "https://show.bilibili.com/platform/detail.html?id=89588"
$ws4.Range("I37").Value = "//i0.hdslb.com/bfs/openplatform/202407/2kN5bTGE1721377069804.png"

# Row 38 <- old row 39 content ("法国姐姐"乔伊丝·乔纳森《小意思》2024巡回演唱会)
$ws4.Range("B38").Value = "'2024-11-17"
$ws4.Range("C38").Value = "广州·“法国姐姐”乔伊丝·乔纳森《小意思》2024巡回演唱会"
$ws4.Range("D38").Value = "东风中路299号 广州中山纪念堂"
$ws4.Range("E38").Value = "2024.11.17 19:30-11.17 21:00"
$ws4.Range("F38").Value = 4
$ws4.Range("G38").Value = 280
$ws4.Range("H38").Value = "https://show.bilibili.com/platform/detail.html?id=91814"
$ws4.Range("I38").Value = "//i2.hdslb.com/bfs/openplatform/202408/bnKPQEEd1725008600562.jpeg"

# Row 39 <- brand-new event (【会员购严选】广州·环形宇宙动漫游戏嘉年华)
$ws4.Range("B39").Value = "'2024-11-23"
$ws4.Range("C39").Value = "【会员购严选】广州·环形宇宙动漫游戏嘉年华"
$ws4.Range("D39").Value = "新港东路630-638号(广交会展旁，地铁新港东站F出口) 南丰国际会展中心"
$ws4.Range("E39").Value = "2024.11.23 09:30-11.24 17:00"
$ws4.Range("F39").Value = 49
$ws4.Range("G39").Value = "不可售"
$ws4.Range("H39").Value = "https://show.bilibili.com/platform/detail.html?id=93064"
$ws4.Range("I39").Value = "//i0.hdslb.com/bfs/openplatform/202409/9M2b0A4e1727589140960.jpeg"
